$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting the latest cryptos data pull (values/volumes refreshed;
# Algorand/Filecoin rows re-ranked and swapped between row 41 and row 42).
$updates = @(
    @{Row=2; D="97.060.75"; E="  +0.75%  "},
    @{Row=3; D="3.690.26"; E="  +0.84%  "},
    @{Row=4; E="  +0.03%  "},
    @{Row=5; D="236.75"; E="  -2.06%  "},
    @{Row=6; D="1.91"; E="  +0.81%  "},
    @{Row=7; D="658.14"; E="  -0.37%  "},
    @{Row=8; E="  +0.21%  "},
    @{Row=9; E="  -2.21%  "},
    @{Row=10; D="1.00"; E="  +0.00%  "},
    @{Row=11; D="3.688.38"; E="  +0.85%  "},
    @{Row=12; D="44.17"; E="  -1.27%  "},
    @{Row=13; D="0.209"; E="  +2.33%  "},
    @{Row=14; D="0.0000304"; E="  +11.90%  "},
    @{Row=15; E="  +1.81%  "},
    @{Row=16; D="4.380.10"; E="  +0.89%  "},
    @{Row=17; D="96.786.69"; E="  +0.70%  "},
    @{Row=18; D="9.11"; E="  +2.23%  "},
    @{Row=19; D="3.704.44"; E="  +0.75%  "},
    @{Row=20; D="12.98"; E="  +2.14%  "},
    @{Row=21; D="18.72"; E="  +2.63%  "},
    @{Row=22; D="0.507"; E="  -4.65%  "},
    @{Row=23; D="519.80"; E="  +0.08%  "},
    @{Row=24; E="  -0.13%  "},
    @{Row=25; E="  +3.19%  "},
    @{Row=26; E="  +0.98%  "},
    @{Row=27; D="0.200"; E="  +19.19%  "},
    @{Row=28; D="101.39"; E="  -0.64%  "},
    @{Row=29; D="13.42"; E="  +3.72%  "},
    @{Row=30; D="12.54"; E="  +1.74%  "},
    @{Row=31; D="3.02"; E="  -0.17%  "},
    @{Row=32; D="0.999"; E="  -0.20%  "},
    @{Row=33; D="0.190"; E="  +2.71%  "},
    @{Row=34; E="  +2.11%  "},
    @{Row=35; D="1.00"; E="  -0.02%  "},
    @{Row=36; D="32.18"; E="  -2.11%  "},
    @{Row=37; D="648.71"; E="  +3.68%  "},
    @{Row=38; D="0.591"; E="  +0.75%  "},
    @{Row=39; D="8.83"; E="  +1.64%  "},
    @{Row=40; E="  +0.01%  "},
    @{Row=41; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="6.82"; E="  +8.92%  "},
    @{Row=42; B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.494"; E="  +9.01%  "},
    @{Row=43; D="2.05"; E="  +5.23%  "},
    @{Row=44; E="  +1.41%  "},
    @{Row=45; D="40.68"; E="  -9.20%  "},
    @{Row=46; D="0.960"; E="  +0.59%  "},
    @{Row=47; D="0.0464"; E="  +1.61%  "},
    @{Row=48; E="  +0.44%  "},
    @{Row=49; D="23.62"; E="  +0.05%  "},
    @{Row=50; D="8.67"; E="  +1.49%  "},
    @{Row=51; E="  -1.13%  "}
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey("B")) { $ws.Cells.Item($row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        $ws.Cells.Item($row, 4).NumberFormat = "@"
        $ws.Cells.Item($row, 4).Value = $u.D
    }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($row, 5).Value = $u.E }
}
